$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Reorder the "Periodo Mora" year labels (E16:E20) from descending (2007..2003)
# to ascending (2003..2007), matching the updated shared-strings order.
$ws.Range("E16").Value = "2003"
$ws.Range("E17").Value = "2004"
$ws.Range("E18").Value = "2005"
$ws.Range("E19").Value = "2006"
$ws.Range("E20").Value = "2007"

# Update "Valor Mora" amounts (G16:G20) to the new base value.
$ws.Range("G16").Value = 828116
$ws.Range("G17").Value = 828116
$ws.Range("G18").Value = 828116
$ws.Range("G19").Value = 828116
$ws.Range("G20").Value = 828116
